$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D; existing D:K data (with formats) shifts right to E:L
$ws.Columns("D").Insert()

# New "most recent period" values for the newly-opened column D.
# Rows not listed here are blank in column D (format-only copy below).
$newValues = @{
    7 = 43465
    8 = 195700
    9 = 102200
    10 = 93500
    12 = 24100
    13 = 0
    14 = -12000
    15 = 2300
    17 = 185600
    18 = 10100
    20 = 400
    21 = 29100
    22 = 100
    23 = 10400
    24 = 2400
    25 = 0
    26 = 8000
    27 = 8000
    28 = 0
    29 = 1800
    30 = 0
    31 = 0
    32 = -400
    33 = 9800
    34 = 0
    35 = 9800
    38 = 43465
    41 = 25400
    42 = 25100
    43 = 37700
    44 = 0
    45 = 3200
    46 = 91400
    47 = 0
    48 = 27400
    49 = 76400
    50 = 0
    51 = 0
    52 = 3700
    53 = 0
    54 = 198900
    57 = 9200
    58 = "NA"
    59 = 23900
    60 = 33100
    61 = 0
    62 = 600
    63 = 0
    64 = 0
    65 = 0
    66 = 33800
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -338600
    73 = 0
    74 = 0
    75 = 0
    76 = 165200
    77 = 0
    80 = 43465
    81 = 9800
    83 = 18600
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 19700
    91 = -16100
    92 = 0
    93 = 0
    94 = -12700
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = -2400
    101 = -100
    102 = 4500
}

$dataRows = 7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102

foreach ($r in $dataRows) {
    $src = $ws.Cells.Item($r, 5)   # column E now holds what used to be in D (same formatting)
    $dst = $ws.Cells.Item($r, 4)   # newly inserted column D
    $src.Copy()
    $dst.PasteSpecial(-4122)       # xlPasteFormats: bring the correct number format/style over
    if ($newValues.ContainsKey($r)) {
        $dst.Value = $newValues[$r]
    }
}

$excel.CutCopyMode = 0

# Column width bookkeeping (characters) to mirror the inserted column & resized neighbours
$ws.Columns.Item(2).ColumnWidth = 26.88671875
$ws.Columns.Item(3).ColumnWidth = 69.109375
for ($c = 4; $c -le 11; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 14.6640625
}
